# Update the "Förändrad" date column (C2:C19) from 2023-10-22 to 2023-10-25
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = Get-Date -Year 2023 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0
$newDate = Get-Date -Year 2023 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value()
    if ($current -eq $oldDate) {
        $cell.Value = $newDate
    }
}
